$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "327.71"
Set-TextValue $ws.Range("E2") "-0.76%"
Set-TextValue $ws.Range("D3") "43.96"
Set-TextValue $ws.Range("E3") "5.56%"
Set-TextValue $ws.Range("D4") "5.477"
Set-TextValue $ws.Range("E4") "-3.85%"
Set-TextValue $ws.Range("D5") "0.08065"
Set-TextValue $ws.Range("E5") "-4.25%"
Set-TextValue $ws.Range("D6") "8.631"
Set-TextValue $ws.Range("E6") "-1.78%"
Set-TextValue $ws.Range("D8") "1.875"
Set-TextValue $ws.Range("E8") "-5.54%"
Set-TextValue $ws.Range("D9") "2.723"
Set-TextValue $ws.Range("E9") "-7.75%"
Set-TextValue $ws.Range("D10") "0.9350"
Set-TextValue $ws.Range("E10") "0.88%"
Set-TextValue $ws.Range("D11") "0.1162"
Set-TextValue $ws.Range("E11") "-8.74%"
Set-TextValue $ws.Range("D12") "0.1891"
Set-TextValue $ws.Range("E12") "-3.71%"
Set-TextValue $ws.Range("D13") "0.09579"
Set-TextValue $ws.Range("E13") "1.77%"
Set-TextValue $ws.Range("D14") "0.04149"
Set-TextValue $ws.Range("E14") "5.14%"
Set-TextValue $ws.Range("D15") "0.1066"
Set-TextValue $ws.Range("E15") "0.23%"
Set-TextValue $ws.Range("D16") "0.001272"
Set-TextValue $ws.Range("E16") "-2.85%"
Set-TextValue $ws.Range("D17") "0.006007"
Set-TextValue $ws.Range("E17") "-1.81%"
Set-TextValue $ws.Range("D18") "3.571"
Set-TextValue $ws.Range("E18") "4.21%"
Set-TextValue $ws.Range("E19") "-0.75%"
Set-TextValue $ws.Range("D20") "8.552"
Set-TextValue $ws.Range("E20") "-4.81%"
Set-TextValue $ws.Range("D21") "0.1366"
Set-TextValue $ws.Range("E21") "0.08%"
Set-TextValue $ws.Range("E22") "3.11%"
Set-TextValue $ws.Range("D23") "0.04330"
Set-TextValue $ws.Range("E23") "-1.92%"
Set-TextValue $ws.Range("D24") "0.001234"
Set-TextValue $ws.Range("E24") "-1.01%"
Set-TextValue $ws.Range("D25") "0.004344"
Set-TextValue $ws.Range("E25") "-1.10%"
Set-TextValue $ws.Range("D26") "0.0001231"
Set-TextValue $ws.Range("E26") "3.24%"
Set-TextValue $ws.Range("D27") "0.0004000"
Set-TextValue $ws.Range("E27") "0.05%"
Set-TextValue $ws.Range("D39") "0.02651"
Set-TextValue $ws.Range("E39") "-6.99%"
Set-TextValue $ws.Range("D40") "0.05440"
Set-TextValue $ws.Range("E40") "-1.36%"
Set-TextValue $ws.Range("D41") "0.01143"
Set-TextValue $ws.Range("E41") "27.20%"
Set-TextValue $ws.Range("D42") "0.007686"
Set-TextValue $ws.Range("E42") "-2.78%"
Set-TextValue $ws.Range("D43") "0.1389"
Set-TextValue $ws.Range("E43") "-3.37%"
Set-TextValue $ws.Range("E44") "1.96%"
Set-TextValue $ws.Range("D45") "0.009660"
Set-TextValue $ws.Range("E45") "-11.91%"
Set-TextValue $ws.Range("D46") "0.00006873"
Set-TextValue $ws.Range("E46") "-5.70%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.05%"
Set-TextValue $ws.Range("D48") "0.003563"
Set-TextValue $ws.Range("E48") "9.65%"
Set-TextValue $ws.Range("D49") "0.002276"
Set-TextValue $ws.Range("E49") "-0.28%"
Set-TextValue $ws.Range("D50") "0.00002106"
Set-TextValue $ws.Range("E50") "0.05%"
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.05%"
